# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '69.578.95'
Set-TextCell "E2" '  +0.48%  '
# Row 3
Set-TextCell "D3" '3.682.14'
Set-TextCell "E3" '  +0.29%  '
# Row 4
Set-TextCell "E4" '  -0.04%  '
# Row 5
Set-TextCell "D5" '667.26'
Set-TextCell "E5" '  -1.11%  '
# Row 6
Set-TextCell "D6" '159.62'
Set-TextCell "E6" '  +0.91%  '
# Row 8
Set-TextCell "D8" '0.500'
Set-TextCell "E8" '  +1.40%  '
# Row 9
Set-TextCell "E9" '  -0.23%  '
# Row 10
Set-TextCell "D10" '7.12'
Set-TextCell "E10" '  +2.54%  '
# Row 11
Set-TextCell "E11" '  +1.36%  '
# Row 12
Set-TextCell "E12" '  +0.95%  '
# Row 13
Set-TextCell "D13" '32.90'
Set-TextCell "E13" '  +1.78%  '
# Row 14
Set-TextCell "D14" '3.665.17'
Set-TextCell "E14" '  -0.32%  '
# Row 15
Set-TextCell "D15" '69.570.91'
# Row 16
Set-TextCell "E16" '  +2.53%  '
# Row 17
Set-TextCell "D17" '16.12'
Set-TextCell "E17" '  +0.48%  '
# Row 18
Set-TextCell "D18" '6.45'
Set-TextCell "E18" '  +0.40%  '
# Row 19
Set-TextCell "D19" '469.90'
Set-TextCell "E19" '  +0.70%  '
# Row 20
Set-TextCell "D20" '9.73'
# Row 21
Set-TextCell "D21" '0.645'
Set-TextCell "E21" '  -0.46%  '
# Row 22
Set-TextCell "D22" '79.66'
Set-TextCell "E22" '  -0.07%  '
# Row 23
Set-TextCell "D23" '3.828.95'
Set-TextCell "E23" '  +0.28%  '
# Row 25
Set-TextCell "D25" '0.0000126'
Set-TextCell "E25" '  +3.46%  '
# Row 26
Set-TextCell "D26" '10.91'
Set-TextCell "E26" '  +0.03%  '
# Row 27
Set-TextCell "D27" '9.03'
Set-TextCell "E27" '  -0.17%  '
# Row 28
Set-TextCell "D28" '2.67'
Set-TextCell "E28" '  -0.25%  '
# Row 29
Set-TextCell "D29" '1.69'
Set-TextCell "E29" '  -2.93%  '
# Row 30
Set-TextCell "D30" '2.00'
Set-TextCell "E30" '  +1.08%  '
# Row 31
Set-TextCell "E31" '  +0.12%  '
# Row 32
Set-TextCell "E32" '  +2.96%  '
# Row 33
Set-TextCell "D33" '26.73'
Set-TextCell "E33" '  -0.48%  '
# Row 34
Set-TextCell "D34" '6.46'
Set-TextCell "E34" '  -2.03%  '
# Row 35
Set-TextCell "D35" '3.684.70'
Set-TextCell "E35" '  +0.53%  '
# Row 36
Set-TextCell "D36" '8.45'
Set-TextCell "E36" '  +3.33%  '
# Row 37
Set-TextCell "D37" '6.07'
Set-TextCell "E37" '  -2.22%  '
# Row 39
Set-TextCell "D39" '2.25'
Set-TextCell "E39" '  +1.43%  '
# Row 40
Set-TextCell "D40" '1.00'
Set-TextCell "E40" '  -0.03%  '
# Row 41
Set-TextCell "D41" '176.95'
Set-TextCell "E41" '  +1.32%  '
# Row 42
Set-TextCell "E42" '  +0.94%  '
# Row 44
Set-TextCell "D44" '47.00'
Set-TextCell "E44" '  -1.13%  '
# Row 45
Set-TextCell "D45" '2.73'
Set-TextCell "E45" '  +1.51%  '
# Row 46
Set-TextCell "D46" '1.28'
Set-TextCell "E46" '  -0.84%  '
# Row 47
Set-TextCell "B47" 'InjectiveProtocol'
Set-TextCell "C47" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell "D47" '27.42'
Set-TextCell "E47" '  -2.14%  '
# Row 48
Set-TextCell "B48" 'FLOKI'
Set-TextCell "C48" 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextCell "D48" '0.000271'
Set-TextCell "E48" '  -2.04%  '
# Row 49
Set-TextCell "B49" 'Cosmos'
Set-TextCell "C49" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell "D49" '7.84'
Set-TextCell "E49" '  +0.82%  '
# Row 50
Set-TextCell "B50" 'SuiNetwork'
Set-TextCell "C50" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell "D50" '1.06'
Set-TextCell "E50" '  -0.67%  '
# Row 51
Set-TextCell "B51" 'TheGraph'
Set-TextCell "C51" 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell "D51" '0.263'
Set-TextCell "E51" '  -0.67%  '
